# Updated symbol list on Mon Dec 12 03:36:45 UTC 2022 with GitHub Actions
#
# This script re-applies the latest price refresh to the "cryptos" sheet:
#  - Column D ("Price") values are text-formatted numeric strings; they are
#    written with NumberFormat "@" (Text) first so Excel does not coerce
#    them into real numbers (which would mangle values like the very small
#    "0.00000000752" into scientific notation), then the cell Style is put
#    back to "Normal" so no stray formatting is left behind.
#  - Rows 42/43 (CEJI / BKEXToken) swapped ranking positions this refresh,
#    so their Coin name, Link and Volume(1h) columns are updated together
#    with their new prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# --- Column D (Price) updates -------------------------------------------------
$priceUpdates = @{
    2  = "283.84"
    3  = "20.71"
    4  = "6.205"
    5  = "0.06171"
    6  = "3.589"
    7  = "6.553"
    8  = "1.506"
    9  = "0.8196"
    10 = "0.01381"
    11 = "0.1637"
    12 = "0.08420"
    13 = "0.03484"
    14 = "0.03214"
    15 = "0.09193"
    16 = "3.716"
    17 = "0.001640"
    18 = "0.04716"
    19 = "0.006435"
    21 = "0.001071"
    22 = "0.0001604"
    23 = "3.843"
    26 = "0.1253"
    40 = "0.04723"
    41 = "0.007215"
    44 = "0.01151"
    45 = "0.00006919"
    46 = "0.00000000752"
    47 = "1.103"
    48 = "0.002825"
    49 = "0.00001906"
    50 = "0.01244"
}

foreach ($row in $priceUpdates.Keys) {
    Set-TextValue $ws.Cells.Item($row, 4) $priceUpdates[$row]
}

# --- Rows 42/43: CEJI and BKEXToken swap ranking positions --------------------
$ws.Cells.Item(42, 2).Value = "BKEXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Cells.Item(42, 4) "0.1101"
$ws.Cells.Item(42, 5).Value = "41BKEXTokenBKK"

$ws.Cells.Item(43, 2).Value = "CEJI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Cells.Item(43, 4) "0.003568"
$ws.Cells.Item(43, 5).Value = "42CEJICEJI"
